# Add the new survey wave (22. 6. 2021) as a new last column on both sheets,
# and bump the "aktualizace" date in the two summary/footer labels from
# "1. 6. 2021" to "28. 6. 2021".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data": percentages table. Existing columns run D..AD (dates from
# "Stav pred epidemii" through "25. 5. 2021"); new column AE holds the
# "22. 6. 2021" wave.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Copy header formatting (bold, centered, bordered) from the previous header
# cell so the new header cell matches the rest of the row, then set its text.
$ws1.Range("AD1").Copy()
$ws1.Range("AE1").PasteSpecial(-4122)
$ws1.Range("AE1").Value = "22. 6. 2021"

$data1 = @{
    2  = 0.09
    3  = 0.08
    4  = 0.1
    5  = 0.07
    6  = 0.08
    7  = 0.13
    8  = 0.09
    9  = 0.17
    10 = 0.1
    11 = 0.07
    12 = 0.12
    13 = 0.07
    14 = 0.2
    15 = 0.1
    16 = 0.07
    17 = 0.13
    18 = 0.11
    19 = 0.05
    20 = 0.1
    21 = 0.05
    22 = 0.06
    23 = 0.18
}

foreach ($row in $data1.Keys) {
    $ws1.Cells.Item($row, 31).Value = $data1[$row]
}

# ---------------------------------------------------------------------------
# Sheet "pocetR": respondent-count table. Existing columns run C..AC; new
# column AD holds the "22. 6. 2021" wave.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

$ws2.Range("AC1").Copy()
$ws2.Range("AD1").PasteSpecial(-4122)
$ws2.Range("AD1").Value = "22. 6. 2021"

$data2 = @{
    2  = 1904
    3  = 937
    4  = 967
    5  = 260
    6  = 677
    7  = 296
    8  = 671
    9  = 161
    10 = 298
    11 = 357
    12 = 342
    13 = 746
    14 = 183
    15 = 377
    16 = 1344
    17 = 194
    18 = 753
    19 = 580
    20 = 249
    21 = 573
    22 = 803
    23 = 528
}

foreach ($row in $data2.Keys) {
    $ws2.Cells.Item($row, 30).Value = $data2[$row]
}

# ---------------------------------------------------------------------------
# Update the footer title strings (row 24, column A on both sheets) to
# reflect the new "update" date, 28. 6. 2021 instead of 1. 6. 2021.
# ---------------------------------------------------------------------------
$oldDate = "1. 6. 2021"
$newDate = "28. 6. 2021"

$title1 = $ws1.Range("A24").Value()
$ws1.Range("A24").Value = $title1.Replace($oldDate, $newDate)

$title2 = $ws2.Range("A24").Value()
$ws2.Range("A24").Value = $title2.Replace($oldDate, $newDate)
